$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new work-log entry in row 2
$ws.Range("A2").Value = "Haley"
$ws.Range("B2").Value = "3 hours"
$ws.Range("C2").Value = "Researching how to use blender, working on textures"

# Match the (slightly) updated column widths recorded in the saved file
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 17.833333333333332
$ws.Columns.Item(3).ColumnWidth = 30.833333333333332

# Update the active selection to A3, matching the saved view state
$ws.Range("A3").Select()
